# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Change B11 on the (single) "Rules" sheet from the text "R40" to the text
# "1". The target cell must keep holding a *text* value ("1"), not the
# number 1, and it must keep its existing cell style/number format
# (General) - i.e. only the stored string changes, nothing else about the
# cell's formatting.
#
# A plain  $ws.Range("B11").Value = "1"  would make Excel auto-detect the
# numeric-looking string and store a real number instead of text, which is
# not what we want here. To force Excel to keep it as text without
# touching B11's own number format, we build the text value on a scratch
# cell (using TEXT() so the result is unambiguously a string), copy it,
# and paste only the *value* (PasteSpecial values) into B11 - this carries
# over the "it's text" nature of the source without changing B11's
# existing style. The scratch cell is removed again afterwards so it
# leaves no trace in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z100")
$scratch.Formula = "=TEXT(1,""0"")"

$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues

$scratch.Delete(-4162)  # xlShiftUp
$excel.CutCopyMode = 0
